$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "unit_price_currency" column (L) data rows were changed from SAR to EUR.
$ws.Range("L2:L153").Value = "EUR"

# Viewport was scrolled down a couple of rows (topLeftCell A128 -> A130)
# while keeping the existing L2:L153 selection untouched.
$excel.ActiveWindow.ScrollRow = 130
